# The workbook is a small "Password" tracker sheet with a header row and
# a free-standing Counter value, plus some leftover blank placeholder
# cells/rows from the previous (Python/openpyxl-generated) version of the
# file. This pass, made by re-opening and editing the sheet in Excel:
#   - bumps the Counter value in F1 from 2 to 3 (one more slot used up)
#   - wipes out the stray blank placeholder cells that openpyxl had left
#     behind (the filler cells I1:R1 under the merged note cell, the
#     lone styled cell A3, and the blank rows 4/5 in B:C) now that they
#     are no longer needed
#   - leaves the cursor/selection on C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the Counter value (F1): this is the actual data change.
$counter = $ws.Range("F1").Value()
$ws.Range("F1").Value = $counter + 1

# Drop the unused blank placeholder cells/rows left over from the old
# version of the sheet.
$ws.Range("A3").Clear()
$ws.Range("B4:C5").Clear()
$ws.Range("I1:R1").ClearContents()
$ws.Range("I1:R1").ClearFormats()

# Leave the selection on C7.
$ws.Range("C7").Select()
